$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple in-place property updates ---
$ws.Range("B3").Value  = "0.1.7"                                   # Version
$ws.Range("B6").Value  = "draft"                                   # Status
$ws.Range("B8").Value  = "2024-08-27T12:23:18-05:00"                # Date
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"  # Contact (publisher)
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"            # Contact (person)

# --- make room for a new "Jurisdiction" row at row 12 ---
# First extend the formatting one row further down (copy row 15's style to the
# brand-new row 16) so the grown table keeps a consistent look, then shift the
# existing Description/Purpose/Copyright/Immutable rows down by one (bottom-up
# so we never clobber a row before it has been copied).
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 15; $r -ge 12; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value()
}

# --- populate the new Jurisdiction row ---
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
